$d = $word.ActiveDocument

function Replace-ParagraphXml($para, $bodyInnerXml) {
    # Replace an entire paragraph's content (incl. its own pPr/runs) with
    # new content described by $bodyInnerXml (one or more <w:p> elements,
    # or <w:p>'s content for a single-paragraph replace). Extends the
    # range to the real end of the document when the paragraph is the very
    # last one in the body (its Range.End excludes the trailing mark there),
    # so the operation does not leave a stray empty paragraph behind.
    $rng = $para.Range
    $docEnd = $d.Content.End
    if ($rng.End -lt $docEnd) {
        $endPos = $rng.End
    } else {
        $endPos = $docEnd
    }
    $full = $d.Range($rng.Start, $endPos)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $full.InsertXML($xml)
}

# ------------------------------------------------------------------
# 1. Insert two brand-new paragraphs at the very start of the body:
#    - a paragraph with the "Name of all cities -> ... -> heatmap..."
#      text (incl. Wingdings arrow symbols) and the _GoBack bookmark
#    - an empty paragraph right after it
# ------------------------------------------------------------------
$firstParaRange = $d.Paragraphs(1).Range
$firstParaRange.Collapse(1)          # wdCollapseStart
$firstParaRange.InsertParagraphBefore()
$firstParaRange.InsertParagraphBefore()

$newPara1Xml = '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Name of all cities </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> call Beijing AQI API </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>for each city</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> store it in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>heatArray</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> heatmap of all cities in selected country.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

Replace-ParagraphXml $d.Paragraphs(1) $newPara1Xml

$emptyPara2Xml = '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
Replace-ParagraphXml $d.Paragraphs(2) $emptyPara2Xml

# ------------------------------------------------------------------
# 2. Remove the _GoBack bookmark from its old location (the empty
#    paragraph that used to hold it, right before "Provide
#    information about their current condition on air quality.")
# ------------------------------------------------------------------
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $next = $d.Paragraphs($i + 1)
    if ($p.Range.Text -eq [string][char]13 -and $next.Range.Text.StartsWith("Provide information about their current condition")) {
        $emptyXml = '<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
        Replace-ParagraphXml $p $emptyXml
        break
    }
}

# ------------------------------------------------------------------
# 3. Add a lastRenderedPageBreak before the "and more particularly
#    the surface concentration..." run near the end of the document.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("and more particularly the surface concentration")) {
        $brokenXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00AF3449"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:lastRenderedPageBreak/><w:t>and more particularly the surface concentration which matters in order to quantify the health impact. Moreover, the tropospheric Ozone is having a diurnal cycle, with pollution peaking in the afternoon when the temperature reaches its maximum, and almost no pollution during the night.</w:t></w:r></w:p>'
        Replace-ParagraphXml $p $brokenXml
        break
    }
}
